$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Sheet1 (o_10): add new column E (evaluator_partial_correctness) ---
$ws1.Range("E1").Value = 'evaluator_partial_correctness'
$ws1.Range("D1").Copy()
$ws1.Range("E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Sheet1 (o_10): update existing row 2 data (prompt/llm_response/evaluator_response changed) ---
$ws1.Range("A2").Value = ' Given is the adjacency matrix for a unweighted undirected graph containing 10 nodes labelled A to J. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   

Consider some examples

Example 1: what is the shortest path from node A to node O?
   A B C D E F G H I J K L M N O
 A 0 1 0 1 0 0 0 0 0 0 0 0 1 0 0
 B 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 1 0 0 0 1 0 0 0 1 0 0 1 0 0 0
 E 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0
 I 0 0 0 1 0 0 0 0 0 1 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0
 L 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 M 1 0 0 0 0 0 0 0 0 0 0 0 0 1 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0

Solution: A -> M -> N -> O
        

Example 2: what is the shortest path from node A to node H?
   A B C D E F G H
 A 0 1 0 0 1 0 1 0
 B 1 0 1 1 0 0 0 0
 C 0 1 0 0 0 0 0 0
 D 0 1 0 0 0 0 0 0
 E 1 0 0 0 0 1 0 0
 F 0 0 0 0 1 0 0 0
 G 1 0 0 0 0 0 0 1
 H 0 0 0 0 0 0 1 0

Solution: A -> G -> H
        

Example 3: what is the shortest path from node A to node O?
   A B C D E F G H I J K L M N O
 A 0 1 0 0 0 1 0 0 0 0 1 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 1 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 F 1 0 0 0 0 0 1 0 0 1 0 0 0 0 0
 G 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0
 J 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 K 1 0 0 0 0 0 0 0 0 0 0 1 1 0 0
 L 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 1 0 0 1 1
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0

Solution: A -> K -> M -> O
        
 Given these examples, answer the following quesiton.

what is the shortest path from node A to node J?

   A B C D E F G H I J
 A 0 1 1 0 0 0 0 0 0 0
 B 1 0 0 0 0 0 0 0 0 0
 C 1 0 0 1 0 1 0 0 0 0
 D 0 0 1 0 1 0 0 0 0 0
 E 0 0 0 1 0 0 0 0 0 0
 F 0 0 1 0 0 0 1 0 0 0
 G 0 0 0 0 0 1 0 1 0 0
 H 0 0 0 0 0 0 1 0 1 0
 I 0 0 0 0 0 0 0 1 0 1
 J 0 0 0 0 0 0 0 0 1 0
    '
$ws1.Range("B2").Value = 'A -> C -> F -> G -> H -> I -> J'
$ws1.Range("C2").Value = 'The shortest path from node A to node J is A -> B -> C -> F -> G -> J.'
$ws1.Range("D2").Value = 'Wrong'
$ws1.Range("E2").Value = 'Output: 3/6'

# --- Add sheet2 (o_20) right after sheet1 ---
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "o_20"
$ws1.Range("A1:E1").Copy($ws2.Range("A1:E1"))
$excel.CutCopyMode = $false
$ws2.Range("A2").Value = ' Given is the adjacency matrix for a unweighted undirected graph containing 20 nodes labelled A to T. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   
Consider some examples
Example 1: what is the shortest path from node A to node Y?
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 1 0 1 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 1 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0
 G 0 0 0 0 0 1 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 1 0 0 0 1 1 0 0 1 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 1 0 0 0 0 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 1 0 0 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0
 W 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 1
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0
Solution: A -> E -> F -> W -> Y
Example 2: what is the shortest path from node A to node R?
   A B C D E F G H I J K L M N O P Q R
 A 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 1 0 0 0 1 0 1 0 0 0 0 0 1 0 0 0 0 0
 E 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 1 0 0 0 1 1 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 1 0 0 1 0 1 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 M 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0
Solution: A -> D -> M -> N -> O -> P -> Q -> R
Example 3: what is the shortest path from node A to node Q?
   A B C D E F G H I J K L M N O P Q
 A 0 1 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0
 B 1 0 1 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 1 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 1 0 1 0 1 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 1 0 0 0 1 1 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0
 M 1 0 0 0 0 0 0 0 0 0 0 0 0 1 1 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 1 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0
Solution: A -> M -> O -> P -> Q
 Given these examples, answer the following quesiton.
what is the shortest path from node A to node T?
   A B C D E F G H I J K L M N O P Q R S T
 A 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 1 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 1 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 1 1 1 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 1 0 0 0 1 0 0 1 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 1 0 1 1 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 P 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0
    '
$ws2.Range("B2").Value = 'A -> P -> Q -> R -> S -> T'
$ws2.Range("C2").Value = 'There is no direct connection between node A and node T in the given adjacency matrix. Therefore, there is no shortest path from node A to node T.'
$ws2.Range("D2").Value = 'Wrong'
$ws2.Range("E2").Value = 'Output: 0/1'

# --- Add sheet3 (o_20_jumbled) right after sheet2 ---
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws2)
$ws3.Name = "o_20_jumbled"
$ws1.Range("A1:E1").Copy($ws3.Range("A1:E1"))
$excel.CutCopyMode = $false
$ws3.Range("A2").Value = ' Given is the adjacency matrix for a unweighted undirected graph containing 20 nodes labelled A to T. The value corresponding to each row M and column N represents whether there is a connection between the two nodes, where 0 means no connection.   
Consider some examples
Example 1: what is the shortest path from node A to node Y?
   A B C D E F G H I J K L M N O P Q R S T U V W X Y
 A 0 1 0 1 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 1 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0
 G 0 0 0 0 0 1 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 1 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 1 0 0 0 1 1 0 0 1 0 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 1 0 0 0 0 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 1 0 1 0 0 0 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0
 T 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0
 U 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 1 0 0 0
 V 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0
 W 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 1
 X 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0
 Y 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0
Solution: A -> E -> F -> W -> Y
Example 2: what is the shortest path from node A to node R?
   A B C D E F G H I J K L M N O P Q R
 A 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 1 0 0 0 1 0 1 0 0 0 0 0 1 0 0 0 0 0
 E 0 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 1 0 0 0 1 1 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 1 0 0 1 0 1 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 M 0 0 0 1 0 0 0 0 0 0 0 0 0 1 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0
Solution: A -> D -> M -> N -> O -> P -> Q -> R
Example 3: what is the shortest path from node A to node Q?
   A B C D E F G H I J K L M N O P Q
 A 0 1 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0
 B 1 0 1 0 0 1 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 1 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 G 0 0 0 0 0 1 0 1 0 1 0 0 0 0 0 0 0
 H 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 1 0 0 0 1 1 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0
 M 1 0 0 0 0 0 0 0 0 0 0 0 0 1 1 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 1 0
 P 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0
Solution: A -> M -> O -> P -> Q
 Given these examples, answer the following quesiton.
what is the shortest path from node A to node T?
   A B C D E F G H I J K L M N O P Q R S T
 A 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 B 1 0 1 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 C 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 D 0 0 1 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 E 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 F 0 1 0 0 0 0 1 1 0 0 0 0 0 0 0 0 0 0 0 1
 G 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
 H 0 0 0 0 0 1 0 0 1 1 0 1 0 0 0 0 0 0 0 0
 I 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0
 J 0 0 0 0 0 0 0 1 0 0 1 0 0 0 0 0 0 0 0 0
 K 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0
 L 0 0 0 0 0 0 0 1 0 0 0 0 1 0 0 1 0 0 0 0
 M 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0 0
 N 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0 0 0 0
 O 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 0 0
 P 0 0 0 0 0 0 0 0 0 0 0 1 0 0 0 0 1 0 0 0
 Q 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0 0
 R 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 1 0
 S 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 0 1 0 0
 T 0 0 0 0 0 1 0 0 0 0 0 0 0 0 0 0 0 0 0 0
    '
$ws3.Range("B2").Value = 'A -> B -> F -> T'
$ws3.Range("C2").Value = 'There is no direct connection between node A and node T in the graph. Therefore, there is no shortest path from node A to node T.'
$ws3.Range("D2").Value = 'Wrong'
$ws3.Range("E2").Value = 'Output: 0/4'

# --- Restore active sheet / selection to sheet1 ---
$ws1.Activate()
$ws1.Range("A1").Select()

Write-Output "edit complete"
